$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for column D price cells so values like "1.00" or
# "46.881.84" are not reinterpreted as numbers (matches source data which
# stores these as literal text/inline strings).
$priceCells = @("D2","D3","D5","D6","D7","D8","D10","D14","D15","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D30","D31","D32","D33","D35","D36","D39","D40","D42","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values scraped by the GitHub Actions cryptos-list job.
$ws.Range("D2").Value = "46.881.84"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "2.337.51"
$ws.Range("E3").Value = "  +4.08%  "
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").Value = "305.97"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "96.75"
$ws.Range("E6").Value = "  +2.74%  "
$ws.Range("D7").Value = "0.576"
$ws.Range("E7").Value = "  +1.13%  "
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.54%  "
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("D10").Value = "35.74"
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("E11").Value = "  +1.07%  "
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "2.694.72"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "2.336.84"
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("E16").Value = "  +4.71%  "
$ws.Range("D17").Value = "0.828"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "46.753.79"
$ws.Range("E18").Value = "  +5.16%  "
$ws.Range("D19").Value = "13.52"
$ws.Range("E19").Value = "  +15.12%  "
$ws.Range("D20").Value = "0.0₃0947"
$ws.Range("E20").Value = "  +1.36%  "
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "67.19"
$ws.Range("E22").Value = "  +2.89%  "
$ws.Range("D23").Value = "245.99"
$ws.Range("E23").Value = "  +3.46%  "
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").Value = "1.99"
$ws.Range("E25").Value = "  +1.21%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "42.47"
$ws.Range("E27").Value = "  +14.69%  "
$ws.Range("D28").Value = "2.26"
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").Value = "20.11"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "5.77"
$ws.Range("E31").Value = "  -0.99%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0816"
$ws.Range("E32").Value = "  +4.77%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "148.40"
$ws.Range("E33").Value = "  -0.18%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("D35").Value = "3.17"
$ws.Range("E35").Value = "  -0.23%  "
$ws.Range("D36").Value = "0.112"
$ws.Range("E36").Value = "  +2.60%  "
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").Value = "4.00"
$ws.Range("E39").Value = "  +6.23%  "
$ws.Range("D40").Value = "0.0314"
$ws.Range("E40").Value = "  +5.81%  "
$ws.Range("E41").Value = "  +1.60%  "
$ws.Range("D42").Value = "13.91"
$ws.Range("E42").Value = "  -7.42%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +11.14%  "
$ws.Range("D45").Value = "1.842.36"
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").Value = "83.45"
$ws.Range("E46").Value = "  +2.38%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "75.25"
$ws.Range("E47").Value = "  +9.25%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.196"
$ws.Range("E48").Value = "  +5.51%  "
$ws.Range("D49").Value = "4.93"
$ws.Range("E49").Value = "  +2.29%  "
$ws.Range("D50").Value = "98.86"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").Value = "55.16"
$ws.Range("E51").Value = "  +2.42%  "
